$d = $word.ActiveDocument

# Locate the (single) field in the document and the paragraph that
# contains it, before we delete the field and lose that context.
$field = $d.Fields.Item(1)
$fieldStart = $field.Code.Start

$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($fieldStart -ge $candidate.Range.Start -and $fieldStart -lt $candidate.Range.End) {
        $targetIndex = $i
    }
}

# Remove the field (fldChar begin/end + instrText runs) entirely,
# leaving an empty paragraph behind in its place.
$field.Delete()

# Re-acquire the (now empty) paragraph and insert plain-text runs
# carrying the same literal characters the field's instrText runs used
# to hold - "{", "m", ":'", "prefix", "\t", "suffix", "'", "}" - plus
# the bookmark that used to sit between "prefix" and "\t".
$targetPara = $d.Paragraphs.Item($targetIndex)
$insertionPoint = $targetPara.Range
$insertionPoint.End = $insertionPoint.Start

$newContent = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:r><w:t>{</w:t></w:r>' + `
    '<w:r><w:t>m</w:t></w:r>' + `
    "<w:r><w:t>:'</w:t></w:r>" + `
    '<w:r><w:t>prefix</w:t></w:r>' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
    '<w:r><w:t>\t</w:t></w:r>' + `
    '<w:r><w:t>suffix</w:t></w:r>' + `
    "<w:r><w:t>'</w:t></w:r>" + `
    '<w:r><w:t xml:space="preserve">}</w:t></w:r>' + `
    '</w:p>'

$null = $insertionPoint.InsertXML($newContent)
